$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 37536
$ws.Range("D2").Value = 54292179
$ws.Range("C3").Value = 90589
$ws.Range("D3").Value = 132809643
$ws.Range("C4").Value = 31055
$ws.Range("D4").Value = 45992694
$ws.Range("C5").Value = 8657
$ws.Range("D5").Value = 12867851
$ws.Range("C6").Value = 1979
$ws.Range("D6").Value = 2941006
$ws.Range("C11").Value = 41065
$ws.Range("D11").Value = 55735425
$ws.Range("C12").Value = 9613
$ws.Range("D12").Value = 13905289
$ws.Range("C13").Value = 25850
$ws.Range("D13").Value = 37912243
$ws.Range("C14").Value = 8292
$ws.Range("D14").Value = 12307263
$ws.Range("C15").Value = 2139
$ws.Range("D15").Value = 3180883
$ws.Range("C16").Value = 415
$ws.Range("D16").Value = 611623
$ws.Range("C17").Value = 32
$ws.Range("D17").Value = 48000
$ws.Range("C19").Value = 10177
$ws.Range("D19").Value = 13481637
$ws.Range("C20").Value = 13332
$ws.Range("D20").Value = 19253251
$ws.Range("C21").Value = 31553
$ws.Range("D21").Value = 46309098
$ws.Range("C22").Value = 10199
$ws.Range("D22").Value = 15161555
$ws.Range("C23").Value = 2623
$ws.Range("D23").Value = 3900655
$ws.Range("C26").Value = 11630
$ws.Range("D26").Value = 15538729
$ws.Range("C27").Value = 7613
$ws.Range("D27").Value = 11030784
$ws.Range("C28").Value = 22390
$ws.Range("D28").Value = 32864772
$ws.Range("C29").Value = 7781
$ws.Range("D29").Value = 11579133
$ws.Range("C33").Value = 8262
$ws.Range("D33").Value = 10919546
$ws.Range("C34").Value = 3215
$ws.Range("D34").Value = 4640115
$ws.Range("C35").Value = 7776
$ws.Range("D35").Value = 11356479
$ws.Range("C36").Value = 3165
$ws.Range("D36").Value = 4690461
$ws.Range("C37").Value = 823
$ws.Range("D37").Value = 1225823
$ws.Range("C40").Value = 2448
$ws.Range("D40").Value = 3308825
$ws.Range("C41").Value = 17145
$ws.Range("D41").Value = 24798453
$ws.Range("C42").Value = 50879
$ws.Range("D42").Value = 74593611
$ws.Range("C43").Value = 18940
$ws.Range("D43").Value = 28133750
$ws.Range("C44").Value = 5584
$ws.Range("D44").Value = 8315978
$ws.Range("C45").Value = 1192
$ws.Range("D45").Value = 1778545
$ws.Range("C49").Value = 16604
$ws.Range("D49").Value = 22116357
$ws.Range("C50").Value = 1989
$ws.Range("D50").Value = 2885982
$ws.Range("C51").Value = 6803
$ws.Range("D51").Value = 10001337
$ws.Range("C52").Value = 2326
$ws.Range("D52").Value = 3473918
$ws.Range("C56").Value = 6754
$ws.Range("D56").Value = 9306960
$ws.Range("C57").Value = 921
$ws.Range("D57").Value = 1351584
$ws.Range("C58").Value = 2299
$ws.Range("D58").Value = 3408317
$ws.Range("C59").Value = 912
$ws.Range("D59").Value = 1357501
$ws.Range("C63").Value = 1357
$ws.Range("D63").Value = 1910356
$ws.Range("C64").Value = 15269
$ws.Range("D64").Value = 22060010
$ws.Range("C65").Value = 44508
$ws.Range("D65").Value = 65137119
$ws.Range("C66").Value = 15644
$ws.Range("D66").Value = 23252332
$ws.Range("C67").Value = 4557
$ws.Range("D67").Value = 6787792
$ws.Range("C72").Value = 15027
$ws.Range("D72").Value = 19819124
$ws.Range("C73").Value = 50933
$ws.Range("D73").Value = 74124636
$ws.Range("C74").Value = 144960
$ws.Range("D74").Value = 213577629
$ws.Range("C75").Value = 63233
$ws.Range("D75").Value = 94225270
$ws.Range("C76").Value = 20190
$ws.Range("D76").Value = 30165817
$ws.Range("C77").Value = 4768
$ws.Range("D77").Value = 7123723
$ws.Range("C78").Value = 260
$ws.Range("D78").Value = 385170
$ws.Range("C84").Value = 50397
$ws.Range("D84").Value = 68594770
$ws.Range("C85").Value = 4555
$ws.Range("D85").Value = 6599441
$ws.Range("C86").Value = 11495
$ws.Range("D86").Value = 16889327
$ws.Range("C87").Value = 3863
$ws.Range("D87").Value = 5757406
$ws.Range("C88").Value = 1340
$ws.Range("D88").Value = 2002489
$ws.Range("C89").Value = 286
$ws.Range("D89").Value = 426512
$ws.Range("C92").Value = 5349
$ws.Range("D92").Value = 7192499
$ws.Range("C93").Value = 1581
$ws.Range("D93").Value = 2276432
$ws.Range("C94").Value = 5110
$ws.Range("D94").Value = 7528429
$ws.Range("C95").Value = 1932
$ws.Range("D95").Value = 2877937
$ws.Range("C96").Value = 685
$ws.Range("D96").Value = 1026460
$ws.Range("C97").Value = 180
$ws.Range("D97").Value = 269113
$ws.Range("C100").Value = 3503
$ws.Range("D100").Value = 4640101
$ws.Range("C101").Value = 592
$ws.Range("D101").Value = 881664
$ws.Range("C102").Value = 346
$ws.Range("D102").Value = 516530
$ws.Range("C104").Value = 43
$ws.Range("D104").Value = 64500
$ws.Range("C106").Value = 10707
$ws.Range("D106").Value = 15539972
$ws.Range("C107").Value = 29086
$ws.Range("D107").Value = 42739514
$ws.Range("C108").Value = 9752
$ws.Range("D108").Value = 14501650
$ws.Range("C109").Value = 2675
$ws.Range("D109").Value = 3988707
$ws.Range("C113").Value = 9737
$ws.Range("D113").Value = 12867953
$ws.Range("C114").Value = 30246
$ws.Range("D114").Value = 43622667
$ws.Range("C115").Value = 65866
$ws.Range("D115").Value = 96402156
$ws.Range("C116").Value = 21281
$ws.Range("D116").Value = 31626332
$ws.Range("C117").Value = 6034
$ws.Range("D117").Value = 8990826
$ws.Range("C118").Value = 1117
$ws.Range("D118").Value = 1669271
$ws.Range("C123").Value = 25697
$ws.Range("D123").Value = 34338270
$ws.Range("C124").Value = 35748
$ws.Range("D124").Value = 51601405
$ws.Range("C125").Value = 76398
$ws.Range("D125").Value = 111730131
$ws.Range("C126").Value = 23737
$ws.Range("D126").Value = 35231389
$ws.Range("C127").Value = 6357
$ws.Range("D127").Value = 9446551
$ws.Range("C128").Value = 1224
$ws.Range("D128").Value = 1820411
$ws.Range("C129").Value = 57
$ws.Range("D129").Value = 83728
$ws.Range("C132").Value = 31575
$ws.Range("D132").Value = 41947538
$ws.Range("C133").Value = 13161
$ws.Range("D133").Value = 19052894
$ws.Range("C134").Value = 32179
$ws.Range("D134").Value = 47266940
$ws.Range("C135").Value = 11437
$ws.Range("D135").Value = 16993542
$ws.Range("C136").Value = 2945
$ws.Range("D136").Value = 4390305
$ws.Range("C137").Value = 493
$ws.Range("D137").Value = 733490
$ws.Range("C140").Value = 10764
$ws.Range("D140").Value = 14360195
$ws.Range("C141").Value = 34825
$ws.Range("D141").Value = 50297313
$ws.Range("C142").Value = 80827
$ws.Range("D142").Value = 118428683
$ws.Range("C143").Value = 24253
$ws.Range("D143").Value = 36037162
$ws.Range("C144").Value = 6359
$ws.Range("D144").Value = 9488067
$ws.Range("C145").Value = 1424
$ws.Range("D145").Value = 2118230
$ws.Range("C146").Value = 80
$ws.Range("D146").Value = 119630
$ws.Range("C148").Value = 28999
$ws.Range("D148").Value = 39147146
